$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '70.694.86'
$ws.Range('E2').Value2 = '  +1.69%  '
$ws.Range('D3').Value2 = '3.630.28'
$ws.Range('E3').Value2 = '  +3.75%  '
$ws.Range('D4').Value2 = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value2 = '  +0.15%  '
$ws.Range('D5').Value2 = "'607.07"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value2 = '  +0.68%  '
$ws.Range('E6').Value2 = '  +2.41%  '
$ws.Range('D7').Value2 = "'0.626"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value2 = '  +0.67%  '
$ws.Range('E8').Value2 = '  +0.11%  '
$ws.Range('D9').Value2 = "'0.221"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value2 = '  +10.25%  '
$ws.Range('E10').Value2 = '  +0.17%  '
$ws.Range('D11').Value2 = "'53.83"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value2 = '  +1.46%  '
$ws.Range('D12').Value2 = "'0.0000305"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value2 = '  +2.03%  '
$ws.Range('D13').Value2 = "'9.57"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value2 = '  +1.18%  '
$ws.Range('D14').Value2 = '4.207.63'
$ws.Range('E14').Value2 = '  +3.69%  '
$ws.Range('D15').Value2 = "'683.90"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value2 = '  +15.17%  '
$ws.Range('B16').Value2 = 'WrappedBTC'
$ws.Range('C16').Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value2 = '70.849.45'
$ws.Range('E16').Value2 = '  +1.61%  '
$ws.Range('B17').Value2 = 'Uniswap'
$ws.Range('C17').Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value2 = "'12.93"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value2 = '  +2.10%  '
$ws.Range('D18').Value2 = '3.653.10'
$ws.Range('E18').Value2 = '  +4.67%  '
$ws.Range('D19').Value2 = "'19.00"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value2 = '  +0.16%  '
$ws.Range('E20').Value2 = '  +0.17%  '
$ws.Range('D21').Value2 = "'0.999"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value2 = "'18.78"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value2 = '  +2.99%  '
$ws.Range('D23').Value2 = "'5.39"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value2 = '  +2.22%  '
$ws.Range('D24').Value2 = "'105.21"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value2 = '  +3.79%  '
$ws.Range('D25').Value2 = "'4.63"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value2 = '  +0.21%  '
$ws.Range('E26').Value2 = '  -4.46%  '
$ws.Range('D27').Value2 = "'10.46"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value2 = '  -2.98%  '
$ws.Range('D28').Value2 = "'9.93"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value2 = '  +4.81%  '
$ws.Range('D29').Value2 = "'34.29"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value2 = '  +3.73%  '
$ws.Range('D30').Value2 = "'4.55"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value2 = '  +6.19%  '
$ws.Range('D31').Value2 = "'7.16"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value2 = '  +2.20%  '
$ws.Range('D32').Value2 = "'12.17"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value2 = '  -1.29%  '
$ws.Range('E33').Value2 = '  +0.73%  '
$ws.Range('D34').Value2 = "'63.29"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value2 = '  +0.34%  '
$ws.Range('B35').Value2 = 'Maker'
$ws.Range('C35').Value2 = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value2 = '3.944.91'
$ws.Range('E35').Value2 = '  +5.61%  '
$ws.Range('B36').Value2 = 'PEPE'
$ws.Range('C36').Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D36').Value2 = '0.0₃0866'
$ws.Range('E36').Value2 = '  +6.74%  '
$ws.Range('D37').Value2 = "'0.999"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E38').Value2 = '  -2.35%  '
$ws.Range('D39').Value2 = "'36.68"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value2 = '  +1.34%  '
$ws.Range('E40').Value2 = '  -0.37%  '
$ws.Range('D41').Value2 = "'500.29"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value2 = '  +1.90%  '
$ws.Range('D42').Value2 = "'3.54"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value2 = '  -2.47%  '
$ws.Range('E43').Value2 = '  +2.30%  '
$ws.Range('E44').Value2 = '  +9.40%  '
$ws.Range('D45').Value2 = "'0.0456"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value2 = '  +1.24%  '
$ws.Range('D46').Value2 = "'3.48"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value2 = '  +5.66%  '
$ws.Range('E47').Value2 = '  +0.61%  '
$ws.Range('D48').Value2 = "'8.66"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value2 = '  +3.20%  '
$ws.Range('E49').Value2 = '  -0.16%  '
$ws.Range('D50').Value2 = "'0.000247"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value2 = '  +1.60%  '
$ws.Range('E51').Value2 = '  +1.74%  '
